# Aufgabenverteilung.xlsx - add "User-Testing", "Evaluierung und Überarbeitung
# Axure-Prototyp", "Abschlusspräsentation" and "Abschlussbericht" rows
# (rows 16-19) to the project task-assignment sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell text, in the same order the original author entered it (keeps the
# shared-string table ordering in line with the authored workbook) --------
$ws.Range("A16").Value = "User-Testing"
$ws.Range("B16").Value = "Professor - Carola, Linda, Simon" + [char]10 + "Sekretariat - Carola, Linda" + [char]10 + "Student - Cuong, Konstantin"

$ws.Range("A18").Value = "Abschlusspräsentation"

$ws.Range("A19").Value = "Abschlussbericht"

$ws.Range("A17").Value = "Evaluierung und Überarbeitung Axure-Prototyp"

$ws.Range("B18").Value = "Carola, Cuong, Konstantin, Linda, Simon (siehe ../05_Abschlusspräsentation/Aufteilung wer macht was.docx)"

# B17 reuses the exact same text already used for A15/B15's "Axure-Prototyping"
# responsible-people cell, so Excel will fold it into the same shared string.
$ws.Range("B17").Value = "Carola, Cuong, Konstantin, Linda, Simon (siehe ../Erarbeiten_der_Gestaltungsloesungen/Prototyping Verteilung"

# --- D16: due date 15/01/2019 (serial 43480), formatted like the other
# date cells in column D (dd/mm/yy-ish "m/d/yyyy" numFmt used throughout). ---
$ws.Range("D16").Value = 43480
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Wrap text for the newly-filled "responsible people" cells, matching
# the wrapping already used throughout column B. ---------------------------
$ws.Range("B16").WrapText = $true
$ws.Range("B17").WrapText = $true
$ws.Range("B18").WrapText = $true

# A17's task-name cell keeps the plain (non-wrapped) built-in style, unlike
# its siblings in column A.
$ws.Range("A17").Style = "Excel Built-in Normal"

# --- Row heights: rows 16-18 grew tall enough to show the wrapped
# multi-person cells (matches row 15's existing height of 90). -------------
$ws.Rows.Item(16).RowHeight = 90
$ws.Rows.Item(17).RowHeight = 90
$ws.Rows.Item(18).RowHeight = 90

# --- View state: scrolled down a bit with B19 as the active cell. ---------
$ws.Range("B19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Applied Aufgabenverteilung row 16-19 updates"
